$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update values in column B ---
$ws.Range("B6").Value = "43H03217"
$ws.Range("B7").Value = "13:52:50 02/07/2024"
$ws.Range("B8").Value = "49 Km/h"
$ws.Range("B14").Value = "Quốc Lộ 1A, X. Tam Anh Nam, H. Núi Thành, Quảng Nam"
$ws.Range("B42").Value = "Tổng số xe"
$ws.Range("B43").Value = "Tổng các trạng thái"

# --- Clear out column C (detail) and related column D duplicate cells ---
$cellsToClear = @(
    "C6","D6",
    "C7","D7",
    "C8","D8",
    "C10","D10",
    "C11",
    "C12",
    "C13",
    "C14","D14",
    "C16",
    "C17",
    "C19",
    "C20",
    "C21",
    "C22",
    "C23",
    "C24",
    "C25",
    "C26","D26",
    "C27","D27",
    "C29","D29",
    "C30","D30",
    "C31","D31",
    "C36",
    "C37",
    "C38",
    "C39",
    "C40",
    "C41",
    "C42",
    "C43",
    "C44",
    "C45",
    "C46"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}
